# Landscaping Data - append rows 436:456 (new daily measurements for 7/11-7/13/2025)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (dates) ------------------------------------------------
# Copy the existing date cell's style (s="1", built-in date number format)
# down across the new rows, then overwrite with the real date values.
$ws.Range("A435").Copy($ws.Range("A436:A456"))

$a = New-Object 'object[,]' 21,1
    $a[0,0]=45849
    $a[1,0]=45849
    $a[2,0]=45849
    $a[3,0]=45849
    $a[4,0]=45849
    $a[5,0]=45849
    $a[6,0]=45849
    $a[7,0]=45850
    $a[8,0]=45850
    $a[9,0]=45850
    $a[10,0]=45850
    $a[11,0]=45850
    $a[12,0]=45850
    $a[13,0]=45850
    $a[14,0]=45851
    $a[15,0]=45851
    $a[16,0]=45851
    $a[17,0]=45851
    $a[18,0]=45851
    $a[19,0]=45851
    $a[20,0]=45851
$ws.Range("A436:A456").Value = $a

# --- Columns B:E (Plant_Type, Plant_Size, Low, High) -----------------
$be = New-Object 'object[,]' 21,4
    $be[0,0]="Flowering"; $be[0,1]="Large"; $be[0,2]=74; $be[0,3]=91
    $be[1,0]="Nonflowering"; $be[1,1]="Medium"; $be[1,2]=74; $be[1,3]=91
    $be[2,0]="Nonflowering"; $be[2,1]="Small"; $be[2,2]=74; $be[2,3]=91
    $be[3,0]="Nonflowering"; $be[3,1]="Medium"; $be[3,2]=74; $be[3,3]=91
    $be[4,0]="Nonflowering"; $be[4,1]="Medium"; $be[4,2]=74; $be[4,3]=91
    $be[5,0]="Nonflowering"; $be[5,1]="Large"; $be[5,2]=74; $be[5,3]=91
    $be[6,0]="Tree"; $be[6,1]="Medium"; $be[6,2]=74; $be[6,3]=91
    $be[7,0]="Flowering"; $be[7,1]="Large"; $be[7,2]=72; $be[7,3]=92
    $be[8,0]="Nonflowering"; $be[8,1]="Medium"; $be[8,2]=72; $be[8,3]=92
    $be[9,0]="Nonflowering"; $be[9,1]="Small"; $be[9,2]=72; $be[9,3]=92
    $be[10,0]="Nonflowering"; $be[10,1]="Medium"; $be[10,2]=72; $be[10,3]=92
    $be[11,0]="Nonflowering"; $be[11,1]="Medium"; $be[11,2]=72; $be[11,3]=92
    $be[12,0]="Nonflowering"; $be[12,1]="Large"; $be[12,2]=72; $be[12,3]=92
    $be[13,0]="Tree"; $be[13,1]="Medium"; $be[13,2]=72; $be[13,3]=92
    $be[14,0]="Flowering"; $be[14,1]="Large"; $be[14,2]=71; $be[14,3]=88
    $be[15,0]="Nonflowering"; $be[15,1]="Medium"; $be[15,2]=71; $be[15,3]=88
    $be[16,0]="Nonflowering"; $be[16,1]="Small"; $be[16,2]=71; $be[16,3]=88
    $be[17,0]="Nonflowering"; $be[17,1]="Medium"; $be[17,2]=71; $be[17,3]=88
    $be[18,0]="Nonflowering"; $be[18,1]="Medium"; $be[18,2]=71; $be[18,3]=88
    $be[19,0]="Nonflowering"; $be[19,1]="Large"; $be[19,2]=71; $be[19,3]=88
    $be[20,0]="Tree"; $be[20,1]="Medium"; $be[20,2]=71; $be[20,3]=88
$ws.Range("B436:E456").Value = $be

# --- Column F (Temp_Diff formula = ABS(Low-High)) ---------------------
for ($r = 436; $r -le 456; $r++) {
    $ws.Range("F$r").Formula = "=ABS(D$r-E$r)"
}

# --- Columns G:T (Rain..Pollen) ----------------------------------------
$gt = New-Object 'object[,]' 21,14
    $gt[0,0]=0; $gt[0,1]=0.1; $gt[0,2]="No"; $gt[0,3]=2; $gt[0,4]="Neutral"; $gt[0,5]=8; $gt[0,6]=0.63; $gt[0,7]=75; $gt[0,8]=30.03; $gt[0,9]=15; $gt[0,10]=0.31; $gt[0,11]=9.9; $gt[0,12]=57; $gt[0,13]=0
    $gt[1,0]=0; $gt[1,1]=0.1; $gt[1,2]="No"; $gt[1,3]=3; $gt[1,4]="Neutral"; $gt[1,5]=8; $gt[1,6]=0.63; $gt[1,7]=75; $gt[1,8]=30.03; $gt[1,9]=15; $gt[1,10]=0.31; $gt[1,11]=9.9; $gt[1,12]=57; $gt[1,13]=0
    $gt[2,0]=0; $gt[2,1]=0; $gt[2,2]="No"; $gt[2,3]=3; $gt[2,4]="Dark"; $gt[2,5]=8; $gt[2,6]=0.63; $gt[2,7]=75; $gt[2,8]=30.03; $gt[2,9]=15; $gt[2,10]=0.31; $gt[2,11]=9.9; $gt[2,12]=57; $gt[2,13]=0
    $gt[3,0]=0; $gt[3,1]=0; $gt[3,2]="No"; $gt[3,3]=3; $gt[3,4]="Bright"; $gt[3,5]=8; $gt[3,6]=0.63; $gt[3,7]=75; $gt[3,8]=30.03; $gt[3,9]=15; $gt[3,10]=0.31; $gt[3,11]=9.9; $gt[3,12]=57; $gt[3,13]=0
    $gt[4,0]=0; $gt[4,1]=0; $gt[4,2]="No"; $gt[4,3]=3; $gt[4,4]="Bright"; $gt[4,5]=8; $gt[4,6]=0.63; $gt[4,7]=75; $gt[4,8]=30.03; $gt[4,9]=15; $gt[4,10]=0.31; $gt[4,11]=9.9; $gt[4,12]=57; $gt[4,13]=0
    $gt[5,0]=0; $gt[5,1]=0.2; $gt[5,2]="No"; $gt[5,3]=4; $gt[5,4]="Bright"; $gt[5,5]=8; $gt[5,6]=0.63; $gt[5,7]=75; $gt[5,8]=30.03; $gt[5,9]=15; $gt[5,10]=0.31; $gt[5,11]=9.9; $gt[5,12]=57; $gt[5,13]=0
    $gt[6,0]=0; $gt[6,1]=0.2; $gt[6,2]="No"; $gt[6,3]=1; $gt[6,4]="Dark"; $gt[6,5]=8; $gt[6,6]=0.63; $gt[6,7]=75; $gt[6,8]=30.03; $gt[6,9]=15; $gt[6,10]=0.31; $gt[6,11]=9.9; $gt[6,12]=57; $gt[6,13]=0
    $gt[7,0]=0.56999999999999995; $gt[7,1]=0.2; $gt[7,2]="No"; $gt[7,3]=2; $gt[7,4]="Neutral"; $gt[7,5]=9; $gt[7,6]=0.48; $gt[7,7]=68; $gt[7,8]=30.03; $gt[7,9]=6; $gt[7,10]=0.31; $gt[7,11]=9.9; $gt[7,12]=57; $gt[7,13]=0
    $gt[8,0]=0.56999999999999995; $gt[8,1]=0.1; $gt[8,2]="No"; $gt[8,3]=3; $gt[8,4]="Neutral"; $gt[8,5]=9; $gt[8,6]=0.48; $gt[8,7]=68; $gt[8,8]=30.03; $gt[8,9]=6; $gt[8,10]=0.31; $gt[8,11]=9.9; $gt[8,12]=57; $gt[8,13]=0
    $gt[9,0]=0.56999999999999995; $gt[9,1]=0.2; $gt[9,2]="No"; $gt[9,3]=3; $gt[9,4]="Bright"; $gt[9,5]=9; $gt[9,6]=0.48; $gt[9,7]=68; $gt[9,8]=30.03; $gt[9,9]=6; $gt[9,10]=0.31; $gt[9,11]=9.9; $gt[9,12]=57; $gt[9,13]=0
    $gt[10,0]=0.56999999999999995; $gt[10,1]=0.25; $gt[10,2]="No"; $gt[10,3]=3; $gt[10,4]="Dark"; $gt[10,5]=9; $gt[10,6]=0.48; $gt[10,7]=68; $gt[10,8]=30.03; $gt[10,9]=6; $gt[10,10]=0.31; $gt[10,11]=9.9; $gt[10,12]=57; $gt[10,13]=0
    $gt[11,0]=0.56999999999999995; $gt[11,1]=0.33333333333333331; $gt[11,2]="No"; $gt[11,3]=3; $gt[11,4]="Bright"; $gt[11,5]=9; $gt[11,6]=0.48; $gt[11,7]=68; $gt[11,8]=30.03; $gt[11,9]=6; $gt[11,10]=0.31; $gt[11,11]=9.9; $gt[11,12]=57; $gt[11,13]=0
    $gt[12,0]=0.56999999999999995; $gt[12,1]=0; $gt[12,2]="No"; $gt[12,3]=4; $gt[12,4]="Bright"; $gt[12,5]=9; $gt[12,6]=0.48; $gt[12,7]=68; $gt[12,8]=30.03; $gt[12,9]=6; $gt[12,10]=0.31; $gt[12,11]=9.9; $gt[12,12]=57; $gt[12,13]=0
    $gt[13,0]=0.56999999999999995; $gt[13,1]=0.75; $gt[13,2]="No"; $gt[13,3]=1; $gt[13,4]="Dark"; $gt[13,5]=9; $gt[13,6]=0.48; $gt[13,7]=68; $gt[13,8]=30.03; $gt[13,9]=6; $gt[13,10]=0.31; $gt[13,11]=9.9; $gt[13,12]=57; $gt[13,13]=0
    $gt[14,0]=0.18; $gt[14,1]=0.2; $gt[14,2]="No"; $gt[14,3]=2; $gt[14,4]="Dark"; $gt[14,5]=7; $gt[14,6]=0.56999999999999995; $gt[14,7]=71; $gt[14,8]=30.03; $gt[14,9]=13; $gt[14,10]=0.52; $gt[14,11]=9.9; $gt[14,12]=54; $gt[14,13]=0
    $gt[15,0]=0.18; $gt[15,1]=0.2; $gt[15,2]="No"; $gt[15,3]=3; $gt[15,4]="Neutral"; $gt[15,5]=7; $gt[15,6]=0.56999999999999995; $gt[15,7]=71; $gt[15,8]=30.03; $gt[15,9]=13; $gt[15,10]=0.52; $gt[15,11]=9.9; $gt[15,12]=54; $gt[15,13]=0
    $gt[16,0]=0.18; $gt[16,1]=0.15; $gt[16,2]="No"; $gt[16,3]=3; $gt[16,4]="Neutral"; $gt[16,5]=7; $gt[16,6]=0.56999999999999995; $gt[16,7]=71; $gt[16,8]=30.03; $gt[16,9]=13; $gt[16,10]=0.52; $gt[16,11]=9.9; $gt[16,12]=54; $gt[16,13]=0
    $gt[17,0]=0.18; $gt[17,1]=0.2; $gt[17,2]="No"; $gt[17,3]=3; $gt[17,4]="Bright"; $gt[17,5]=7; $gt[17,6]=0.56999999999999995; $gt[17,7]=71; $gt[17,8]=30.03; $gt[17,9]=13; $gt[17,10]=0.52; $gt[17,11]=9.9; $gt[17,12]=54; $gt[17,13]=0
    $gt[18,0]=0.18; $gt[18,1]=0.1; $gt[18,2]="No"; $gt[18,3]=3; $gt[18,4]="Bright"; $gt[18,5]=7; $gt[18,6]=0.56999999999999995; $gt[18,7]=71; $gt[18,8]=30.03; $gt[18,9]=13; $gt[18,10]=0.52; $gt[18,11]=9.9; $gt[18,12]=54; $gt[18,13]=0
    $gt[19,0]=0.18; $gt[19,1]=0.4; $gt[19,2]="No"; $gt[19,3]=4; $gt[19,4]="Dark"; $gt[19,5]=7; $gt[19,6]=0.56999999999999995; $gt[19,7]=71; $gt[19,8]=30.03; $gt[19,9]=13; $gt[19,10]=0.52; $gt[19,11]=9.9; $gt[19,12]=54; $gt[19,13]=0
    $gt[20,0]=0.18; $gt[20,1]=0.95; $gt[20,2]="No"; $gt[20,3]=1; $gt[20,4]="Bright"; $gt[20,5]=7; $gt[20,6]=0.56999999999999995; $gt[20,7]=71; $gt[20,8]=30.03; $gt[20,9]=13; $gt[20,10]=0.52; $gt[20,11]=9.9; $gt[20,12]=54; $gt[20,13]=0
$ws.Range("G436:T456").Value = $gt

# Row 447's Rain column (H) was entered in the source workbook as the
# formula =1/3 rather than a typed decimal literal - reproduce that.
$ws.Range("H447").Formula = "=1/3"

# --- View state: match the author's final scroll/selection -----------
$ws.Range("T449:T456").Select()
$excel.ActiveWindow.ScrollRow = 432
$excel.ActiveWindow.ScrollColumn = 1
